$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stock names appended to column A (rows 20-24)
$newColA = @(
    "Asian Pain Share Price",
    "wipro share price",
    "lti mindtree share price",
    "sonata software share price",
    "Honeywell share price"
)

for ($i = 0; $i -lt $newColA.Length; $i++) {
    $row = 20 + $i
    $ws.Cells.Item($row, 1).Value = $newColA[$i]
}

# Full column B (rows 1-24): the daily change values for each stock.
# Values that Excel would otherwise auto-convert to numbers (those
# starting with a literal "+") are entered with a leading apostrophe so
# they stay text, matching the "+0.080" style price-change strings.
$colB = @(
    "20/12/2022",
    "−34.70",
    "−147.85",
    "'+0.080",
    "−14.90",
    "−58.55",
    "−4.45",
    "−0.80",
    "'+11.95",
    "−3.05",
    "−7.90",
    "−1.40",
    "−9.65",
    "−30.10",
    "−8.50",
    "−7.25",
    "−57.85",
    "−14.95",
    "'+16.15",
    "−14.05",
    "−2.45",
    "−36.05",
    "−14.65",
    "−0.33"
)

for ($i = 0; $i -lt $colB.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

$ws.Columns.Item(2).AutoFit()

$ws.Range("D8").Select()
